$wb = $excel.ActiveWorkbook

# --- Settings sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Settings")

# Remove the old "OrchestratorQueueName / ProcessABCQueue / queue description"
# row entirely (row 2). Everything below shifts up by one.
$ws.Rows(2).Delete() | Out-Null

# Fill in the new label column for the three URL rows (rows 5-7 after the shift).
$ws.Range("A5").Value = "BoysSourceUrl"
$ws.Range("A6").Value = "GirlsSourceUrl"
$ws.Range("A7").Value = "UnicornNameUrl"

# Add the hyperlinks (Unicorn first, then Girls, then Boys - matches the order
# they were created in when this was originally edited, which drives the
# r:id numbering of the relationships).
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.rpasamples.com/unicornname", "", "", "https://www.rpasamples.com/unicornname") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.gov.pl/web/cyfryzacja/najpopularniejsze-imiona-dla-dziewczynek-2018-ranking-ogolnopolski", "", "", "https://www.gov.pl/web/cyfryzacja/najpopularniejsze-imiona-dla-dziewczynek-2018-ranking-ogolnopolski") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.gov.pl/web/cyfryzacja/najpopularniejsze-imiona-dla-chlopcow-2018-ranking-ogolnopolski", "", "", "https://www.gov.pl/web/cyfryzacja/najpopularniejsze-imiona-dla-chlopcow-2018-ranking-ogolnopolski") | Out-Null

# New "NamesAmount" setting (row 9).
$ws.Range("A9").Value = "NamesAmount"
$ws.Range("B9").Value = 10

# --- Assets sheet -----------------------------------------------------------
# Leave a trailing selection on A25 before moving away from this sheet.
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Select() | Out-Null
$wsAssets.Range("A25").Select() | Out-Null

# --- finish back on Settings, with A9 selected ------------------------------
$ws.Select() | Out-Null
$ws.Range("A9").Select() | Out-Null
